# Add an "_id" column (C) with the Mongo document id for the seeded event,
# matching the "Converted seed data to json" commit: the sheet now carries
# an _id column so the JSON-converted seed can be matched back to its
# Mongo document when the events collection is cleared/reseeded.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header + value for column C.
$ws.Range("C1").Value = "_id"
$ws.Range("C2").Value = "5a847edee5847831acb269a4"

# Size the new column the way Excel would after typing/pasting the id.
$ws.Columns.Item(3).ColumnWidth = 25.5

# Leave the new column selected (entire column), as in the authored file.
$ws.Columns.Item(3).Select() | Out-Null
